# data pegawai, pencarian, sebagian history, ganti password
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Daftar Pegawai")

# New employee rows appended after row 167 (No. Urut 164) -> rows 168-170
$data = @(
    @{ No = 165; Id = 165; Nip = 52635235235;  Nama = "Burhanu Sultan Ramadan"; Tmp = "Darit"; Tgl = "2001-08-09"; Agama = "Islam"; JK = "Laki-laki"; Gol = "B+"; Nikah = "Belum Menikah"; Status = "Aktif"; Pensiun = "2059-08-09" },
    @{ No = 166; Id = 166; Nip = 123456;        Nama = "Jawhead";                 Tmp = "Darit"; Tgl = "2001-08-09"; Agama = "Islam"; JK = "Laki-laki"; Gol = "A+"; Nikah = "Sudah Menikah"; Status = "Aktif"; Pensiun = "2059-08-09" },
    @{ No = 167; Id = 167; Nip = 1234567890;    Nama = "Badang";                  Tmp = "Darit"; Tgl = "2001-08-09"; Agama = "Islam"; JK = "Laki-laki"; Gol = "A-"; Nikah = "Sudah Menikah"; Status = "Aktif"; Pensiun = "2059-08-09" }
)

$row = 168
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec.No
    $ws.Cells.Item($row, 2).Value = $rec.Id
    $ws.Cells.Item($row, 3).Value = $rec.Nip
    $ws.Cells.Item($row, 4).Value = $rec.Nama
    $ws.Cells.Item($row, 5).Value = $rec.Tmp

    $ws.Cells.Item($row, 6).NumberFormat = "@"
    $ws.Cells.Item($row, 6).Value = $rec.Tgl

    $ws.Cells.Item($row, 7).Value = $rec.Agama
    $ws.Cells.Item($row, 8).Value = $rec.JK
    $ws.Cells.Item($row, 9).Value = $rec.Gol
    $ws.Cells.Item($row, 10).Value = $rec.Nikah
    $ws.Cells.Item($row, 11).Value = $rec.Status

    # Touch L:Q (Alamat, Telp, Email, Jabatan, Pendidikan, Unit Kerja) so the
    # cells exist in the sheet (matching the template's blank-cell layout)
    # without minting a new cell style.
    $ws.Cells.Item($row, 12).WrapText = $false
    $ws.Cells.Item($row, 13).WrapText = $false
    $ws.Cells.Item($row, 14).WrapText = $false
    $ws.Cells.Item($row, 15).WrapText = $false
    $ws.Cells.Item($row, 16).WrapText = $false
    $ws.Cells.Item($row, 17).WrapText = $false

    $ws.Cells.Item($row, 18).NumberFormat = "@"
    $ws.Cells.Item($row, 18).Value = $rec.Pensiun

    $row++
}
